$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 2.67
$ws.Range("AD3").Value = 15
$ws.Range("AH3").Value = 25
$ws.Range("AI3").Value = 13.5
$ws.Range("AL3").Value = 50
$ws.Range("AS3").Value = 23

# Row 4 updates
$ws.Range("G4").Value = 1.9
$ws.Range("I4").Value = 3.7
$ws.Range("J4").Value = 2.63
$ws.Range("L4").Value = 4.75
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("AC4").Value = 6
$ws.Range("AD4").Value = 8
$ws.Range("AF4").Value = 15
$ws.Range("AI4").Value = 8
$ws.Range("AJ4").Value = 6.5
$ws.Range("AN4").Value = 9.5
$ws.Range("AO4").Value = 19

# Row 5 updates
$ws.Range("G5").Value = 1.65
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 2.3
$ws.Range("K5").Value = 2.1
$ws.Range("L5").Value = 5.5
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 8.5
$ws.Range("S5").Value = 2.1
$ws.Range("T5").Value = 1.7
$ws.Range("AA5").Value = 2
$ws.Range("AB5").Value = 1.73
$ws.Range("AC5").Value = 6
$ws.Range("AD5").Value = 7
$ws.Range("AE5").Value = 8.5
$ws.Range("AF5").Value = 12
$ws.Range("AH5").Value = 29
$ws.Range("AJ5").Value = 7
$ws.Range("AK5").Value = 19
$ws.Range("AN5").Value = 12
$ws.Range("AO5").Value = 26
$ws.Range("AP5").Value = 17
$ws.Range("AQ5").Value = 51
$ws.Range("AR5").Value = 41

$wb.Save()
